$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header values for columns L, M, N
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Copy style from an existing header cell (K1) to the new header cells
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# New data values for columns L, M, N per row
$ws.Range("L2").Value = 91.56965423913746
$ws.Range("M2").Value = 255984
$ws.Range("N2").Value = 317.5980148883374

$ws.Range("L3").Value = 82.2979860710347
$ws.Range("M3").Value = 7569
$ws.Range("N3").Value = 315.375

$ws.Range("L4").Value = 91.21586592230445
$ws.Range("M4").Value = 185734
$ws.Range("N4").Value = 140.9210925644917

$ws.Range("L5").Value = 79.29010140385255
$ws.Range("M5").Value = 17912
$ws.Range("N5").Value = 275.5692307692308

$ws.Range("L6").Value = 19.6055125364595
$ws.Range("M6").Value = 2083
$ws.Range("N6").Value = 13.88666666666667

$ws.Range("L7").Value = 18.78940113071737
$ws.Range("M7").Value = 125
$ws.Range("N7").Value = 62.5
